$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7721.8
$ws.Range("I40").Value = 5073.6665
$ws.Range("K40").Value = 5073.6665
$ws.Range("M40").Value = -4898.6665
$ws.Range("H86").Value = 3443.6667
$ws.Range("I86").Value = 3441.25
$ws.Range("J86").Value = 3445.6
$ws.Range("K86").Value = 3441.25
$ws.Range("L86").Value = 3445.6
$ws.Range("M86").Value = -2318.25
$ws.Range("N86").Value = -5691.6
$ws.Range("H89").Value = 3443.6667
$ws.Range("I89").Value = 3441.25
$ws.Range("J89").Value = 3445.6
$ws.Range("K89").Value = 17206.25
$ws.Range("L89").Value = 17228
$ws.Range("M89").Value = -11590.25
$ws.Range("N89").Value = -28460
$ws.Range("H98").Value = 1133.4
$ws.Range("I98").Value = 1170.5555
$ws.Range("J98").Value = 799
$ws.Range("K98").Value = 1170.5555
$ws.Range("L98").Value = 799
$ws.Range("M98").Value = 327.4445000000001
$ws.Range("N98").Value = -3795
$ws.Range("H122").Value = 1133.4
$ws.Range("I122").Value = 1170.5555
$ws.Range("J122").Value = 799
$ws.Range("K122").Value = 3511.6665
$ws.Range("L122").Value = 2397
$ws.Range("M122").Value = -1061.6665
$ws.Range("N122").Value = -7297
$ws.Range("H132").Value = 3429.111
$ws.Range("J132").Value = 3995.3333
$ws.Range("L132").Value = 11985.9999
$ws.Range("N132").Value = -17045.9999
$ws.Range("H135").Value = 502.625
$ws.Range("I135").Value = 591
$ws.Range("J135").Value = 308.2
$ws.Range("K135").Value = 5319
$ws.Range("L135").Value = 2773.8
$ws.Range("M135").Value = -2784
$ws.Range("N135").Value = -7843.799999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 825.2222
$ws.Range("J2").Value = 978.3
$ws.Range("L2").Value = 978.3
$ws.Range("N2").Value = -1204.3
$ws.Range("H61").Value = 4158
$ws.Range("I61").Value = 3890.6
$ws.Range("K61").Value = 3890.6
$ws.Range("M61").Value = -3678.6
$ws.Range("H88").Value = 1252
$ws.Range("I88").Value = 1170
$ws.Range("K88").Value = 1170
$ws.Range("M88").Value = -764
$ws.Range("H91").Value = 1252
$ws.Range("I91").Value = 1170
$ws.Range("K91").Value = 1170
$ws.Range("M91").Value = 234
$ws.Range("H110").Value = 619.7143
$ws.Range("I110").Value = 389.83334
$ws.Range("K110").Value = 389.83334
$ws.Range("M110").Value = 1655.16666
$ws.Range("H116").Value = 825.2222
$ws.Range("J116").Value = 978.3
$ws.Range("L116").Value = 978.3
$ws.Range("N116").Value = -5566.3
$ws.Range("H122").Value = 1520.9474
$ws.Range("I122").Value = 1080.0714
$ws.Range("K122").Value = 3240.2142
$ws.Range("M122").Value = -790.2142000000003
$ws.Range("H132").Value = 1847.5714
$ws.Range("I132").Value = 1625.2
$ws.Range("J132").Value = 2403.5
$ws.Range("K132").Value = 4875.6
$ws.Range("L132").Value = 7210.5
$ws.Range("M132").Value = -2345.6
$ws.Range("N132").Value = -12270.5
$ws.Range("H136").Value = 4158
$ws.Range("I136").Value = 3890.6
$ws.Range("K136").Value = 11671.8
$ws.Range("M136").Value = -9121.799999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 825.2222
$ws.Range("J3").Value = 978.3
$ws.Range("L3").Value = 978.3
$ws.Range("N3").Value = -1206.3
$ws.Range("H20").Value = 1213.5834
$ws.Range("I20").Value = 1084.7778
$ws.Range("J20").Value = 1600
$ws.Range("K20").Value = 1084.7778
$ws.Range("L20").Value = 1600
$ws.Range("M20").Value = -837.7778000000001
$ws.Range("N20").Value = -2094
$ws.Range("H122").Value = 150000
$ws.Range("J122").Value = 150000
$ws.Range("L122").Value = 150000
$ws.Range("N122").Value = -159800
$ws.Range("H125").Value = 44000
$ws.Range("J125").Value = 44000
$ws.Range("L125").Value = 44000
$ws.Range("N125").Value = -53840
$ws.Range("H134").Value = 5365
$ws.Range("I134").Value = 5486.6665
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 16459.9995
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -13924.9995
$ws.Range("N134").Value = -20070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5649.971
$ws.Range("J31").Value = 7096.4
$ws.Range("L31").Value = 7096.4
$ws.Range("N31").Value = -7686.4
$ws.Range("H34").Value = 5649.971
$ws.Range("J34").Value = 7096.4
$ws.Range("L34").Value = 7096.4
$ws.Range("N34").Value = -7500.4
$ws.Range("H132").Value = 2448.111
$ws.Range("I132").Value = 2129.125
$ws.Range("K132").Value = 6387.375
$ws.Range("M132").Value = -3857.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H113").Value = 750.0714
$ws.Range("I113").Value = 546.9091
$ws.Range("K113").Value = 1640.7273
$ws.Range("M113").Value = 529.2727
$ws.Range("H132").Value = 70
$ws.Range("I132").Value = 70
$ws.Range("K132").Value = 630
$ws.Range("M132").Value = 1900

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2999
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 2999
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 4749
$ws.Range("I80").Value = 4749
$ws.Range("K80").Value = 4749
$ws.Range("M80").Value = -3751
$ws.Range("H83").Value = 4749
$ws.Range("I83").Value = 4749
$ws.Range("K83").Value = 23745
$ws.Range("M83").Value = -18753
$ws.Range("H102").Value = 1267.0714
$ws.Range("J102").Value = 495
$ws.Range("L102").Value = 495
$ws.Range("N102").Value = -3739
$ws.Range("H111").Value = 10000
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 7400
$ws.Range("I113").Value = 3500
$ws.Range("K113").Value = 3500
$ws.Range("M113").Value = -1330
$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -31884
$ws.Range("H126").Value = 2972.875
$ws.Range("I126").Value = 2299
$ws.Range("K126").Value = 6897
$ws.Range("M126").Value = -4427
$ws.Range("H132").Value = 4199
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6650
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 7820
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 7820
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -8196
$ws.Range("H68").Value = 7624.875
$ws.Range("I68").Value = 5249.75
$ws.Range("K68").Value = 5249.75
$ws.Range("M68").Value = -4500.75
$ws.Range("H71").Value = 7624.875
$ws.Range("I71").Value = 5249.75
$ws.Range("K71").Value = 26248.75
$ws.Range("M71").Value = -22504.75
$ws.Range("H93").Value = 3277.5715
$ws.Range("I93").Value = 2490.5
$ws.Range("K93").Value = 2490.5
$ws.Range("M93").Value = -1242.5
$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7974.625
$ws.Range("I81").Value = 5883
$ws.Range("J81").Value = 14249.5
$ws.Range("K81").Value = 11766
$ws.Range("L81").Value = 28499
$ws.Range("M81").Value = -10705
$ws.Range("N81").Value = -30621
$ws.Range("H84").Value = 7974.625
$ws.Range("I84").Value = 5883
$ws.Range("J84").Value = 14249.5
$ws.Range("K84").Value = 58830
$ws.Range("L84").Value = 142495
$ws.Range("M84").Value = -53526
$ws.Range("N84").Value = -153103
$ws.Range("H132").Value = 2096.3845
$ws.Range("I132").Value = 1773.1111
$ws.Range("K132").Value = 5319.3333
$ws.Range("M132").Value = -2789.3333
